$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Test Status" value from PASS to FAIL
$ws.Range("Z2").Value = "FAIL"

# Highlight the failing status in red (was a light green/teal, indexed 42 -> red, indexed 10)
$ws.Range("Z2").Interior.ColorIndex = 3

# Scroll testing for elements: bring the view back to show column I and
# move the active selection to I2 (was scrolled to show column S with S2 selected)
[void]$ws.Range("I2").Select()
